# Apply the F-column classification swaps as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value  = "Outras compras"
$ws.Range("F11").Value = "Tecnologia da informação"
$ws.Range("F27").Value = "Outras compras"
$ws.Range("F28").Value = "Outras compras"
$ws.Range("F50").Value = "Tecnologia da informação"
$ws.Range("F71").Value = "Outras compras"
$ws.Range("F73").Value = "Tecnologia da informação"
$ws.Range("F84").Value = "Outras compras"
$ws.Range("F90").Value = "Outras compras"
